$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999467352230126
$ws.Range("E2").Value = 0.9999467352230126

# Row 3
$ws.Range("D3").Value = 0.9998918169644847
$ws.Range("E3").Value = 0.9998918169644847

# Row 4
$ws.Range("D4").Value = 0.2806800994624609
$ws.Range("E4").Value = 0.2806800994624609

# Row 5
$ws.Range("D5").Value = 0.5954755138490253
$ws.Range("E5").Value = 0.5954755138490253

# Row 6 (values in 1E-06 range; write as division to avoid bare exponent literal syntax)
$ws.Range("D6").Value = 7.128983212896906 / 1000000
$ws.Range("E6").Value = 7.128983212896906 / 1000000

# Row 7
$ws.Range("D7").Value = 0.6393701536904405
$ws.Range("E7").Value = 0.3606298463095595

# Row 8
$ws.Range("D8").Value = 0.0005819994459206901
$ws.Range("E8").Value = 0.9994180005540793

# Row 9
$ws.Range("D9").Value = 0.995404981262248
$ws.Range("E9").Value = 0.004595018737752032

# Row 11
$ws.Range("D11").Value = 0.9999962901056528
$ws.Range("E11").Value = 3.709894347170284 / 1000000
$ws.Range("F11").Value = 2.810734272003174
